$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview!G2 - Latest HO Xliff Generate Date (also shared text with de-de!H2)
$overview.Range("G2").Value = "2016-09-04 15:10:00"

# zh-cn!H2 - Correspond Handoff Datetime
$zhcn.Range("H2").Value = "2016-09-04 15:09:55"

# zh-cn!K2 - Correspond Handback DateTime
$zhcn.Range("K2").Value = "2016-09-04 15:10:21"

# de-de!H2 - Correspond Handoff Datetime (shares string with Overview!G2 new value)
$dede.Range("H2").Value = "2016-09-04 15:10:00"

# de-de!K2 - Correspond Handback DateTime
$dede.Range("K2").Value = "2016-09-04 15:10:28"
